$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2 and 3: Sales Document No. changes to 2025002968 ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025002968"
$ws.Range("A2").ClearFormats()
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025002968"
$ws.Range("A3").ClearFormats()

# --- Add new row 4 (duplicate of row 3, SD Test Case) ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025002968"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "13"
$ws.Range("B4").ClearFormats()
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0"
$ws.Range("D4").ClearFormats()
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Y101"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "Sales Order(Prd)"
$ws.Range("H4").ClearFormats()
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "2400"
$ws.Range("I4").ClearFormats()
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "YEA"
$ws.Range("J4").ClearFormats()
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "30"
$ws.Range("K4").ClearFormats()
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = "Inter-company"
$ws.Range("L4").ClearFormats()
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "GN40"
$ws.Range("M4").ClearFormats()
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = "Overseas Int.Company"
$ws.Range("N4").ClearFormats()
$ws.Range("S4").NumberFormat = "@"
$ws.Range("S4").Value = "5800- TEST 16"
$ws.Range("S4").ClearFormats()
$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value = "12-09-2025"
$ws.Range("U4").ClearFormats()
$ws.Range("V4").NumberFormat = "@"
$ws.Range("V4").Value = "31-10-2024"
$ws.Range("V4").ClearFormats()
$ws.Range("X4").NumberFormat = "@"
$ws.Range("X4").Value = "Y12"
$ws.Range("X4").ClearFormats()
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "Y12 Order Intake with PO"
$ws.Range("Y4").ClearFormats()
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "0.00"
$ws.Range("Z4").ClearFormats()
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "EUR"
$ws.Range("AA4").ClearFormats()
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "31-10-2024"
$ws.Range("AB4").ClearFormats()
$ws.Range("AK4").NumberFormat = "@"
$ws.Range("AK4").Value = "Y2U00027"
$ws.Range("AK4").ClearFormats()
$ws.Range("AM4").NumberFormat = "@"
$ws.Range("AM4").Value = "Y2U00027"
$ws.Range("AM4").ClearFormats()
$ws.Range("AO4").NumberFormat = "@"
$ws.Range("AO4").Value = "200132028"
$ws.Range("AO4").ClearFormats()
$ws.Range("AP4").NumberFormat = "@"
$ws.Range("AP4").Value = "Yokogawa India (Pty) Ltd."
$ws.Range("AP4").ClearFormats()
$ws.Range("AQ4").NumberFormat = "@"
$ws.Range("AQ4").Value = "0"
$ws.Range("AQ4").ClearFormats()
$ws.Range("BA4").NumberFormat = "@"
$ws.Range("BA4").Value = "End User"
$ws.Range("BA4").ClearFormats()
$ws.Range("BB4").NumberFormat = "@"
$ws.Range("BB4").Value = "Name: End User"
$ws.Range("BB4").ClearFormats()
$ws.Range("BG4").NumberFormat = "@"
$ws.Range("BG4").Value = "2"
$ws.Range("BG4").ClearFormats()
$ws.Range("BI4").NumberFormat = "@"
$ws.Range("BI4").Value = "0"
$ws.Range("BI4").ClearFormats()
$ws.Range("BM4").NumberFormat = "@"
$ws.Range("BM4").Value = "0"
$ws.Range("BM4").ClearFormats()
$ws.Range("BO4").NumberFormat = "@"
$ws.Range("BO4").Value = "2"
$ws.Range("BO4").ClearFormats()
$ws.Range("BW4").NumberFormat = "@"
$ws.Range("BW4").Value = "SAPP"
$ws.Range("BW4").ClearFormats()
$ws.Range("BX4").NumberFormat = "@"
$ws.Range("BX4").Value = "SO: Approved"
$ws.Range("BX4").ClearFormats()
$ws.Range("CA4").NumberFormat = "@"
$ws.Range("CA4").Value = "0"
$ws.Range("CA4").ClearFormats()
$ws.Range("CE4").NumberFormat = "@"
$ws.Range("CE4").Value = "Customer Group"
$ws.Range("CE4").ClearFormats()
$ws.Range("CF4").NumberFormat = "@"
$ws.Range("CF4").Value = "Text: Customer Group"
$ws.Range("CF4").ClearFormats()
$ws.Range("CN4").NumberFormat = "@"
$ws.Range("CN4").Value = "30059222"
$ws.Range("CN4").ClearFormats()
$ws.Range("CO4").NumberFormat = "@"
$ws.Range("CO4").Value = "Rohit Bharad"
$ws.Range("CO4").ClearFormats()
$ws.Range("CP4").NumberFormat = "@"
$ws.Range("CP4").Value = "12-09-2024"
$ws.Range("CP4").ClearFormats()
$ws.Range("CT4").NumberFormat = "@"
$ws.Range("CT4").Value = "10"
$ws.Range("CT4").ClearFormats()
$ws.Range("CW4").NumberFormat = "@"
$ws.Range("CW4").Value = "F3XD64_F000000001"
$ws.Range("CW4").ClearFormats()
$ws.Range("DB4").NumberFormat = "@"
$ws.Range("DB4").Value = "10.000"
$ws.Range("DB4").ClearFormats()
$ws.Range("DD4").NumberFormat = "@"
$ws.Range("DD4").Value = "130"
$ws.Range("DD4").ClearFormats()
$ws.Range("DE4").NumberFormat = "@"
$ws.Range("DE4").Value = "661.00"
$ws.Range("DE4").ClearFormats()
$ws.Range("DF4").NumberFormat = "@"
$ws.Range("DF4").Value = "100.00"
$ws.Range("DF4").ClearFormats()
$ws.Range("DG4").NumberFormat = "@"
$ws.Range("DG4").Value = "0"
$ws.Range("DG4").ClearFormats()
$ws.Range("DI4").NumberFormat = "@"
$ws.Range("DI4").Value = "0.00"
$ws.Range("DI4").ClearFormats()
$ws.Range("DJ4").NumberFormat = "@"
$ws.Range("DJ4").Value = "0.00"
$ws.Range("DJ4").ClearFormats()
$ws.Range("DL4").NumberFormat = "@"
$ws.Range("DL4").Value = "0"
$ws.Range("DL4").ClearFormats()
$ws.Range("DM4").NumberFormat = "@"
$ws.Range("DM4").Value = "%"
$ws.Range("DM4").ClearFormats()
$ws.Range("DP4").NumberFormat = "@"
$ws.Range("DP4").Value = "0.00"
$ws.Range("DP4").ClearFormats()
$ws.Range("DS4").NumberFormat = "@"
$ws.Range("DS4").Value = "0.00"
$ws.Range("DS4").ClearFormats()
$ws.Range("DV4").NumberFormat = "@"
$ws.Range("DV4").Value = "0.00"
$ws.Range("DV4").ClearFormats()
$ws.Range("DY4").NumberFormat = "@"
$ws.Range("DY4").Value = "0.00"
$ws.Range("DY4").ClearFormats()
$ws.Range("EB4").NumberFormat = "@"
$ws.Range("EB4").Value = "0.00"
$ws.Range("EB4").ClearFormats()
$ws.Range("EJ4").NumberFormat = "@"
$ws.Range("EJ4").Value = "0.00"
$ws.Range("EJ4").ClearFormats()
$ws.Range("EM4").NumberFormat = "@"
$ws.Range("EM4").Value = "0.00"
$ws.Range("EM4").ClearFormats()
$ws.Range("EP4").NumberFormat = "@"
$ws.Range("EP4").Value = "0.00"
$ws.Range("EP4").ClearFormats()
$ws.Range("ES4").NumberFormat = "@"
$ws.Range("ES4").Value = "0.00"
$ws.Range("ES4").ClearFormats()
$ws.Range("EV4").NumberFormat = "@"
$ws.Range("EV4").Value = "0.00"
$ws.Range("EV4").ClearFormats()
$ws.Range("EW4").NumberFormat = "@"
$ws.Range("EW4").Value = "0"
$ws.Range("EW4").ClearFormats()
$ws.Range("EY4").NumberFormat = "@"
$ws.Range("EY4").Value = "0.00"
$ws.Range("EY4").ClearFormats()
$ws.Range("EZ4").NumberFormat = "@"
$ws.Range("EZ4").Value = "0"
$ws.Range("EZ4").ClearFormats()
$ws.Range("FB4").NumberFormat = "@"
$ws.Range("FB4").Value = "0.00"
$ws.Range("FB4").ClearFormats()
$ws.Range("FC4").NumberFormat = "@"
$ws.Range("FC4").Value = "0"
$ws.Range("FC4").ClearFormats()
$ws.Range("FE4").NumberFormat = "@"
$ws.Range("FE4").Value = "0.00"
$ws.Range("FE4").ClearFormats()
$ws.Range("FF4").NumberFormat = "@"
$ws.Range("FF4").Value = "0"
$ws.Range("FF4").ClearFormats()
$ws.Range("FH4").NumberFormat = "@"
$ws.Range("FH4").Value = "0.00"
$ws.Range("FH4").ClearFormats()
$ws.Range("FI4").NumberFormat = "@"
$ws.Range("FI4").Value = "0"
$ws.Range("FI4").ClearFormats()
$ws.Range("FK4").NumberFormat = "@"
$ws.Range("FK4").Value = "0.00"
$ws.Range("FK4").ClearFormats()
$ws.Range("FL4").NumberFormat = "@"
$ws.Range("FL4").Value = "0.00"
$ws.Range("FL4").ClearFormats()
$ws.Range("FM4").NumberFormat = "@"
$ws.Range("FM4").Value = "0.00"
$ws.Range("FM4").ClearFormats()
$ws.Range("FN4").NumberFormat = "@"
$ws.Range("FN4").Value = "0.00"
$ws.Range("FN4").ClearFormats()
$ws.Range("FO4").NumberFormat = "@"
$ws.Range("FO4").Value = "0.00"
$ws.Range("FO4").ClearFormats()
$ws.Range("FP4").NumberFormat = "@"
$ws.Range("FP4").Value = "0.00"
$ws.Range("FP4").ClearFormats()
$ws.Range("FQ4").NumberFormat = "@"
$ws.Range("FQ4").Value = "0.00"
$ws.Range("FQ4").ClearFormats()
$ws.Range("FR4").NumberFormat = "@"
$ws.Range("FR4").Value = "0.000"
$ws.Range("FR4").ClearFormats()
$ws.Range("FS4").NumberFormat = "@"
$ws.Range("FS4").Value = "%"
$ws.Range("FS4").ClearFormats()
$ws.Range("FT4").NumberFormat = "@"
$ws.Range("FT4").Value = "0.00"
$ws.Range("FT4").ClearFormats()
$ws.Range("HF4").NumberFormat = "@"
$ws.Range("HF4").Value = "987654321"
$ws.Range("HF4").ClearFormats()
$ws.Range("HK4").NumberFormat = "@"
$ws.Range("HK4").Value = "2002"
$ws.Range("HK4").ClearFormats()
$ws.Range("HM4").NumberFormat = "@"
$ws.Range("HM4").Value = "FCA"
$ws.Range("HM4").ClearFormats()
$ws.Range("HN4").NumberFormat = "@"
$ws.Range("HN4").Value = "0"
$ws.Range("HN4").ClearFormats()
$ws.Range("IG4").NumberFormat = "@"
$ws.Range("IG4").Value = "P"
$ws.Range("IG4").ClearFormats()
$ws.Range("IH4").NumberFormat = "@"
$ws.Range("IH4").Value = "N"
$ws.Range("IH4").ClearFormats()
$ws.Range("II4").NumberFormat = "@"
$ws.Range("II4").Value = "P"
$ws.Range("II4").ClearFormats()
$ws.Range("IJ4").NumberFormat = "@"
$ws.Range("IJ4").Value = "N"
$ws.Range("IJ4").ClearFormats()
$ws.Range("IK4").NumberFormat = "@"
$ws.Range("IK4").Value = "N"
$ws.Range("IK4").ClearFormats()
$ws.Range("IL4").NumberFormat = "@"
$ws.Range("IL4").Value = "1.000"
$ws.Range("IL4").ClearFormats()
$ws.Range("IN4").NumberFormat = "@"
$ws.Range("IN4").Value = "N"
$ws.Range("IN4").ClearFormats()
$ws.Range("IO4").NumberFormat = "@"
$ws.Range("IO4").Value = "N"
$ws.Range("IO4").ClearFormats()
$ws.Range("IP4").NumberFormat = "@"
$ws.Range("IP4").Value = "661.00"
$ws.Range("IP4").ClearFormats()
$ws.Range("IQ4").NumberFormat = "@"
$ws.Range("IQ4").Value = "USD"
$ws.Range("IQ4").ClearFormats()
$ws.Range("IR4").NumberFormat = "@"
$ws.Range("IR4").Value = "Not Performed"
$ws.Range("IR4").ClearFormats()
$ws.Range("JT4").NumberFormat = "@"
$ws.Range("JT4").Value = "0.00"
$ws.Range("JT4").ClearFormats()

# --- GB4: preserve fill style (s=7) from GB3, value unchanged from row 3 ---
$ws.Range("GB3").Copy()
$ws.Range("GB4").PasteSpecial(-4122)
$ws.Range("GB4").Value = "F3XD64-3F/K2/CT"

# --- HJ4: preserve date style from HJ3, set date value ---
$ws.Range("HJ3").Copy()
$ws.Range("HJ4").PasteSpecial(-4122)
$ws.Range("HJ4").Value = 46021

$excel.CutCopyMode = 0
